# The workbook has a single data sheet ("données14"). A logic problem in
# the data generation meant some rows under-estimated column C relative to
# column A, so both columns are corrected row-by-row (column B is untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $aValue, $cValue) {
    $ws.Range("A$row").Value = $aValue
    $ws.Range("C$row").Value = $cValue
}

Set-Row 18  5.4399999999999995   96
Set-Row 19  27.01                89
Set-Row 20  27.439999999999998   93
Set-Row 24  24.12                79
Set-Row 25  23.189999999999998   96
Set-Row 30  62.639999999999993   90
Set-Row 32  42.99                99
Set-Row 35  4.74                 100
Set-Row 36  6.5600000000000005   66
Set-Row 38  32.51                93
Set-Row 39  14.82                99
Set-Row 44  42.19                88
Set-Row 46  43.72                87
Set-Row 47  26.88                96
Set-Row 48  12.83                100
Set-Row 50  7.6499999999999995   61
Set-Row 55  25.490000000000002   96
Set-Row 56  14.469999999999999   96
Set-Row 57  33.54                99
Set-Row 58  21.279999999999998   89
Set-Row 59  6.69                 83
Set-Row 60  8.0299999999999994   96
Set-Row 61  30.14                98
Set-Row 64  28.62                92
Set-Row 71  6.4399999999999995   100
Set-Row 73  10.209999999999999   99
Set-Row 76  23.880000000000003   89
Set-Row 77  9.39                 68
Set-Row 78  13                   100
Set-Row 80  19.12                92
Set-Row 83  7.6899999999999995   99
Set-Row 84  18.72                90
Set-Row 87  11.84                96
Set-Row 88  87.22999999999999    99
Set-Row 89  12.790000000000001   100
Set-Row 92  45.04                100
Set-Row 94  22.759999999999998   86
Set-Row 95  21.51                101
Set-Row 99  45.79                85
Set-Row 100 2.59                 97
